$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "realtionship" -> "relationship"
$ws.Range("A5").Value = "more discussion of the relationship of the optimal policy function to the TR"

# Append two new to-do items at the bottom of the list
$ws.Range("A19").Value = "see loss for RE-optimal TR coefficients under learning"
$ws.Range("A20").Value = "what is truly my message? Once I know, what would I want a paper to do to convince me of this message?"

# Move the selection down to the new last row, as in the authored edit
[void]$ws.Range("A21").Select()
